# Mise à jour de l'application
# Add a new attendance-tracking column (CQ) for the next training session,
# mirroring the previous column (CP)'s formatting, and carrying forward
# each player's last recorded status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header date row) ------------------------------------------
# Set the new session date first (so dependent formulas recalc), then
# copy CP1's number format / style onto the new cell.
$ws.Range("CQ1").Value = 46008
$ws.Range("CP1").Copy()
$ws.Range("CQ1").PasteSpecial(-4122)

# --- Data rows that carry the same attendance mark as column CP -------
$rowsWithValue = @(2,3,4,5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29)
foreach ($r in $rowsWithValue) {
    $src = $ws.Range("CP$r")
    $dst = $ws.Range("CQ$r")
    # Write the value first so the dependency graph / formulas downstream
    # (COUNTA / COUNTIF over K:VQ) pick up the change on recalculation...
    $dst.Value = $src.Value2
    # ...then copy the source cell's formatting only, reusing its style.
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# --- Row 21: no attendance recorded yet for this player in the new ----
# column, but the cell still exists with the same style as CP21.
$ws.Range("CP21").Copy($ws.Range("CQ21"))

# Row 12 intentionally has no cell in column CQ (matches the source data).

$excel.CutCopyMode = $false

# Update the active selection to match the edit's end state.
$ws.Range("CS22").Select()
